# DPE_CARLA_INES.xlsx edit
# - Lower-case all the short variable-code strings in column F / H (rows 19-36)
# - Fix the GJ -> kcal conversion factor (4.2 -> 4.184) (row 19, col H)
# - Combine "KD,\nKM" into "kd;km" (semicolon, no spacing) (row 28, col F)
# - Add "mna" indication to MK -> "mk;mna" (row 36, col F)
# - Row 28 no longer needs the extra wrapped-text height
# - Update the active selection / scrolled position of the sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# Re-assign cell values in the exact order the strings were first introduced
# so that the shared-string table comes out in the same order as the
# authoritative file.
$ws.Range("F19").Value = "gj"
$ws.Range("F20").Value = "zk"
$ws.Range("F21").Value = "ze"
$ws.Range("F22").Value = "zf"
$ws.Range("F23").Value = "za"
$ws.Range("F24").Value = "zb"
$ws.Range("F25").Value = "fs"
$ws.Range("F26").Value = "fu"
$ws.Range("F27").Value = "fp"
$ws.Range("F31").Value = "kmt"
$ws.Range("F32").Value = "kmf"
$ws.Range("F35").Value = "mna"
$ws.Range("H28").Value = "kd+km"
$ws.Range("H19").Value = "gj/4.184"
$ws.Range("F28").Value = "kd;km"
$ws.Range("F36").Value = "mk;mna"

# Row 28 previously needed ht="30" to show the wrapped "KD,\nKM" label; the
# new single-line "kd;km" value fits on one line, so let Excel shrink the
# row back down to the default height.
$ws.Rows.Item(28).AutoFit()

# Move the view: scrolled position + active cell/selection.
[void]$ws.Range("F42").Select()
